$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text (some new values like "1.006" would
# otherwise be auto-coerced to numbers by Excel); reset style afterwards
# so cells don't pick up a stray style index.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$priceData = @{
    2 = "28.150.56"
    3 = "1.780.50"
    4 = "1.006"
    5 = "338.91"
    6 = "1.000"
    7 = "0.3811"
    8 = "0.3441"
    9 = "47.13"
    10 = "1.150"
    11 = "0.07384"
    12 = "23.16"
    13 = "1.003"
    14 = "6.429"
    15 = "7.271"
    16 = "1.779.48"
    17 = "0.00001072"
    18 = "0.06660"
    19 = "82.30"
    20 = "1.000"
    21 = "17.50"
    22 = "6.451"
    23 = "28.186.79"
    24 = "12.04"
    25 = "2.379"
    26 = "1.438"
    27 = "20.75"
    28 = "2.413"
    29 = "153.69"
    30 = "136.62"
    31 = "1.982.04"
    32 = "6.121"
    33 = "3.939"
    34 = "0.08872"
    35 = "12.76"
    36 = "0.02425"
    37 = "0.6831"
    38 = "5.304"
    39 = "0.06347"
    40 = "0.2166"
    41 = "1.242"
    42 = "1.500"
    43 = "8.266"
    44 = "14.23"
    45 = "1.000"
    46 = "0.6273"
    47 = "3.877"
    48 = "132.83"
    49 = "2.085"
    50 = "0.07449"
    51 = "1.203"
}

$volumeData = @{
    2 = "  +3.61%  "
    3 = "  -0.20%  "
    4 = "  +0.22%  "
    5 = "  +0.42%  "
    6 = "  -0.09%  "
    7 = "  -1.28%  "
    8 = "  +0.26%  "
    9 = "  -1.63%  "
    10 = "  -3.38%  "
    11 = "  -0.97%  "
    12 = "  +6.78%  "
    13 = "  +0.20%  "
    14 = "  -0.23%  "
    15 = "  +2.02%  "
    16 = "  -0.08%  "
    17 = "  -1.81%  "
    18 = "  +0.02%  "
    19 = "  -1.30%  "
    20 = "  -0.04%  "
    21 = "  -0.32%  "
    22 = "  -1.00%  "
    23 = "  +3.77%  "
    24 = "  -2.55%  "
    25 = "  +0.55%  "
    26 = "  -0.53%  "
    27 = "  -1.89%  "
    28 = "  -3.36%  "
    29 = "  -1.93%  "
    30 = "  +1.76%  "
    31 = "  -0.05%  "
    32 = "  +2.08%  "
    33 = "  -0.83%  "
    34 = "  +2.05%  "
    35 = "  -1.66%  "
    36 = "  +3.41%  "
    37 = "  +0.12%  "
    38 = "  -1.82%  "
    39 = "  +0.12%  "
    40 = "  -1.16%  "
    41 = "  +0.47%  "
    42 = "  -7.41%  "
    43 = "  -2.16%  "
    44 = "  -0.07%  "
    45 = "  +0.01%  "
    46 = "  -2.11%  "
    47 = "  +0.47%  "
    48 = "  +0.96%  "
    49 = "  -3.90%  "
    50 = "  +4.74%  "
    51 = "  +7.75%  "
}

foreach ($r in 2..51) {
    $ws.Cells.Item($r, 4).Value = $priceData[$r]
    $ws.Cells.Item($r, 5).Value = $volumeData[$r]
}

# Restore default (style-less) formatting on column D now that the
# values are locked in as text.
$dRange.Style = "Normal"
